# Update countries & provincias Spain
# Refresh COVID case counters for a handful of countries and the
# "last updated" timestamp banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 11:37"

# Row 18: Banglades
$ws.Range("B18").Value = 356767
$ws.Range("C18").Value = 1383
$ws.Range("D18").Value = 267024
$ws.Range("E18").Value = 84650
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 5093

# Row 25: Alemania
$ws.Range("B25").Value = 281503
$ws.Range("C25").Value = 158
$ws.Range("E25").Value = 22483
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 9520

# Row 26: Indonesia
$ws.Range("B26").Value = 266845
$ws.Range("C26").Value = 4823
$ws.Range("D26").Value = 196196
$ws.Range("E26").Value = 60431
$ws.Range("G26").Value = 113
$ws.Range("H26").Value = 10218

# Row 27: Israel
$ws.Range("B27").Value = 215273
$ws.Range("C27").Value = 3158
$ws.Range("D27").Value = 152837
$ws.Range("E27").Value = 61031
$ws.Range("G27").Value = 27
$ws.Range("H27").Value = 1405

# Row 47: Polonia
$ws.Range("B47").Value = 84396
$ws.Range("C47").Value = 1587
$ws.Range("D47").Value = 66740
$ws.Range("E47").Value = 15264
$ws.Range("G47").Value = 23
$ws.Range("H47").Value = 2392

# Row 67: Austria
$ws.Range("B67").Value = 41500
$ws.Range("C67").Value = 684
$ws.Range("D67").Value = 32301
$ws.Range("E67").Value = 8413
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 786

# Row 69: Afganistan
$ws.Range("B69").Value = 39186
$ws.Range("C69").Value = 16
$ws.Range("E69").Value = 5116

# Row 77: El Salvador
$ws.Range("B77").Value = 28415
$ws.Range("C77").Value = 214
$ws.Range("D77").Value = 22643
$ws.Range("E77").Value = 4949

# Row 122: Hong Kong
$ws.Range("B122").Value = 5059
$ws.Range("C122").Value = 2
$ws.Range("D122").Value = 4765
$ws.Range("E122").Value = 190

# Row 124: Congo
$ws.Range("B124").Value = 5007
$ws.Range("C124").Value = 192
$ws.Range("D124").Value = 3353
$ws.Range("E124").Value = 1509
$ws.Range("H124").Value = 145

# Row 125: Eslovenia
$ws.Range("B125").Value = 5005
$ws.Range("D125").Value = 3887
$ws.Range("E125").Value = 1029
$ws.Range("H125").Value = 89

# Row 133: Lituania
$ws.Range("B133").Value = 4184
$ws.Range("C133").Value = 114
$ws.Range("D133").Value = 2298
$ws.Range("E133").Value = 1797

# Row 157: Principado de Andorra
$ws.Range("B157").Value = 1758
$ws.Range("C157").Value = 52
$ws.Range("D157").Value = 1072
$ws.Range("E157").Value = 664
$ws.Range("H157").Value = 22

# Row 158: Togo
$ws.Range("B158").Value = 1753
$ws.Range("D158").Value = 1203
$ws.Range("E158").Value = 497
$ws.Range("H158").Value = 53

# Row 159: Belice
$ws.Range("B159").Value = 1707
$ws.Range("D159").Value = 1307
$ws.Range("E159").Value = 356
$ws.Range("H159").Value = 44

# Row 161: Letonia
$ws.Range("D161").Value = 1282
$ws.Range("E161").Value = 307

# Row 175: Tanzania
$ws.Range("B175").Value = 510
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 480
$ws.Range("E175").Value = 23
$ws.Range("H175").Value = 7

# Row 176: Taiwan
$ws.Range("D176").Value = 183
$ws.Range("E176").Value = 305
$ws.Range("H176").Value = 21

# Row 215: Montserrat
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# Row 216: Islas Malvinas
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
